# "Extended airport model added"
# Rework the departure-time / gate schedule on Sheet1:
#   - Gates now rotate 1,2,3,4 within each departure-time group instead of
#     alternating 1,2
#   - Several flights are consolidated onto shared departure times
#   - The last flight's departure time moves earlier (22:20 -> 17:00)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (DepartureTime, Gate) values for rows 2..13 - column A (Destination)
# is untouched.
$ws.Range("B2").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 13 -Minute 20 -Second 0
$ws.Range("C2").Value  = 1

$ws.Range("B3").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 13 -Minute 20 -Second 0
$ws.Range("C3").Value  = 2

$ws.Range("B4").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 13 -Minute 20 -Second 0
$ws.Range("C4").Value  = 3

$ws.Range("B5").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 13 -Minute 20 -Second 0
$ws.Range("C5").Value  = 4

$ws.Range("B6").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 15 -Second 0
$ws.Range("C6").Value  = 1

$ws.Range("B7").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 15 -Second 0
$ws.Range("C7").Value  = 2

$ws.Range("B8").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 15 -Second 0
$ws.Range("C8").Value  = 3

$ws.Range("B9").Value  = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 15 -Second 0
$ws.Range("C9").Value  = 4

$ws.Range("B10").Value = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 45 -Second 0
$ws.Range("C10").Value = 1

$ws.Range("B11").Value = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 45 -Second 0
$ws.Range("C11").Value = 2

$ws.Range("B12").Value = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 14 -Minute 45 -Second 0
$ws.Range("C12").Value = 3

$ws.Range("B13").Value = Get-Date -Year 2015 -Month 12 -Day 21 -Hour 17 -Minute 0 -Second 0
$ws.Range("C13").Value = 4

# Minor column-width tweak (A: 17.5546875 -> 17.5 chars-equivalent)
$ws.Columns.Item(1).ColumnWidth = 16.64
$ws.Columns.Item(2).ColumnWidth = 18.33

# View state: zoom in and move the active selection to B7
$ws.Activate()
$excel.ActiveWindow.Zoom = 191
$ws.Range("B7").Select()
